# Auto update: 2025-12-05 13:19:50
# Updates the latest scan values on Sheet1 of the 국장_조선_분석 workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (HDKSOE) - price/RSI/return refresh plus recomputed scores
$ws.Range("D2").Value = 426500
$ws.Range("E2").Value = 46.4
$ws.Range("F2").Value = 4.02
$ws.Range("K2").Value = 57.7
$ws.Range("N2").Value = 54.83846622768671

# Row 3 (HD HYUNDAI MIPO) - recomputed final score / macro score
$ws.Range("K3").Value = 56.5
$ws.Range("N3").Value = 54.83846622768671

# Row 4 (Hanwha Ocean) - recomputed final score / macro score
$ws.Range("K4").Value = 51.5
$ws.Range("N4").Value = 54.83846622768671

# Row 5 (SamsungHvyInd) - recomputed final score / macro score
$ws.Range("K5").Value = 47.7
$ws.Range("N5").Value = 54.83846622768671
